$d = $word.ActiveDocument

# --- 1. Move the _GoBack bookmark so it becomes a zero-length bookmark
#        right after its start (i.e. its end now sits inside the first
#        paragraph, immediately after the start, instead of after the
#        whole paragraph).
$bm = $d.Bookmarks("_GoBack")
$startPos = $bm.Start
$bm.End = $startPos

# --- 2. Clean up leftover spell-check / grammar-check markers that were
#        splitting runs of plain text; merge them back into single runs
#        by replacing the split text with the unsplit text.
$d.Content.Find.Execute(
    "the comments at the top of platform.h, and the constants in platform.c.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the comments at the top of platform.h, and the constants in platform.c.",
    2) | Out-Null

$d.Content.Find.Execute(
    "The analog co-processor chip has its Tx connected to the Rx of both the base board and KitProg2 and vice versa so that it can communicate via UART to either the base board or to the KitProg2. Therefore, the base board cannot communicate over UART to the KitProg2 since the Tx/Rx lines would be reversed.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The analog co-processor chip has its Tx connected to the Rx of both the base board and KitProg2 and vice versa so that it can communicate via UART to either the base board or to the KitProg2. Therefore, the base board cannot communicate over UART to the KitProg2 since the Tx/Rx lines would be reversed.",
    2) | Out-Null

$d.Content.Find.Execute(
    "To open the workspace in PSoC Creator, double-click on the workspace (cywrk) file. Note, you must have PSoC Creator 4.0 or later installed to open the project.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "To open the workspace in PSoC Creator, double-click on the workspace (cywrk) file. Note, you must have PSoC Creator 4.0 or later installed to open the project.",
    2) | Out-Null

# --- 3. Update the cached NUMPAGES field result in the footer (total
#        page count) from 6 to 1.
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute("6", $true, $false, $false, $false, $false,
                            $true, 1, $false, "1", 2) | Out-Null
